$wb = $excel.ActiveWorkbook

# "MonsterTypes" sheet (sheet2.xml) gets two new columns: goldMult and xpMult,
# each filled with a multiplier of 1 for every monster type row.
$ws = $wb.Worksheets.Item("MonsterTypes")

$ws.Cells.Item(1, 9).Value = "goldMult"
$ws.Cells.Item(1, 10).Value = "xpMult"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = 1
}

# Selection moves to J10 on the MonsterTypes sheet, and this sheet becomes
# the active / selected tab (instead of MonsterRarity).
$ws.Range("J10").Select()
$ws.Activate()

$wb.Save()
